$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 0.8
$ws.Range("A3").NumberFormat = "d-mmm"

$ws.Range("A6").Select()
